# Day_12.xlsx edit: add a second HLOOKUP reference table (vlookup/hlookup practice)
# and re-point the "vishal mega mart" / "big basket" rows at it via HLOOKUP formulas.
# Also tightens the Sheet1 VLOOKUP shared-formula span (H7:H12 -> H7:H11).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Sheet1: the VLOOKUP shared formula in column H was anchored H7:H12, but H12
# was never actually part of that formula (it holds the literal label
# "vlookup"). Re-enter the formula only across the real H7:H11 span.
# ---------------------------------------------------------------------------
$ws1.Range("H7:H11").Formula = "=VLOOKUP(C7,K`$5:M`$11,3,FALSE)"

# ---------------------------------------------------------------------------
# Sheet2: add the second lookup table (shop_name/mode/employee/outlets/sales/
# address/customer_count) in C17:I19, matching the layout/style of the
# existing C13:I15 table.
# ---------------------------------------------------------------------------
$ws2.Range("C13:I15").Copy($ws2.Range("C17"))

$ws2.Range("C17").Value = "shop_name"
$ws2.Range("D17").Value = "mode"
$ws2.Range("E17").Value = "employee"
$ws2.Range("F17").Value = "outlets"
$ws2.Range("G17").Value = "sales"
$ws2.Range("H17").Value = "address"
$ws2.Range("I17").Value = "customer_count"

$ws2.Range("C18").Value = "vishal mega mart"
$ws2.Range("D18").Value = "offline"
$ws2.Range("E18").Value = 100
$ws2.Range("F18").Value = 2
$ws2.Range("G18").Value = 10000
$ws2.Range("H18").Value = "bhande plot"
$ws2.Range("I18").Value = 1000

$ws2.Range("C19").Value = "big basket"
$ws2.Range("D19").Value = "online"
$ws2.Range("E19").Value = 7000
$ws2.Range("F19").Value = 3
$ws2.Range("G19").Value = 150000
$ws2.Range("H19").Value = "wardhman "
$ws2.Range("I19").Value = 12000

# ---------------------------------------------------------------------------
# Row 5 ("vishal mega mart") previously held literal values; replace with
# HLOOKUP formulas against the new C17:I19 table (lookup row 2 = mode).
# ---------------------------------------------------------------------------
$ws2.Range("D5").Formula = "=HLOOKUP(D3,`$C`$17:`$I`$19,2,FALSE)"
$ws2.Range("E5:I5").Formula = "=HLOOKUP(E3,`$C`$17:`$I`$19,2,FALSE)"

# ---------------------------------------------------------------------------
# Row 8 ("big basket") previously held literal values; replace with HLOOKUP
# formulas against the new C17:I19 table (lookup row 3 = employee).
# ---------------------------------------------------------------------------
$ws2.Range("D8:E8").Formula = "=HLOOKUP(D3,`$C`$17:`$I`$19,3,FALSE)"
$ws2.Range("F8").Formula = "=HLOOKUP(F3,`$C`$17:`$I`$19,3,FALSE)"
$ws2.Range("G8:I8").Formula = "=HLOOKUP(G3,`$C`$17:`$I`$19,3,FALSE)"

# ---------------------------------------------------------------------------
# View cosmetics: the sheet grew to C1:I19, scrolled down one row, zoomed to
# 110%, with the selection left on E11.
# ---------------------------------------------------------------------------
$ws2.Activate()
$excel.ActiveWindow.Zoom = 110
$ws2.Range("E11").Select()
